$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- E38 gains a quantity*price-style value (167) ---
$ws.Range("E38").Value = 167

# --- Row 39: new component (solder) ---
$ws.Range("D30").Copy() | Out-Null
$ws.Range("D39").PasteSpecial(-4122) | Out-Null
$ws.Range("D39").Value = 1
$ws.Range("E39").Value = 350
$url39 = "https://dratek.cz/arduino/7545-cinova-pajka-bezolovnata-100g-sn99-3cu0-7.html "
$ws.Range("F39").Value = $url39
$ws.Hyperlinks.Add($ws.Range("F39"), $url39.Trim()) | Out-Null
$ws.Range("F30").Copy() | Out-Null
$ws.Range("F39").PasteSpecial(-4122) | Out-Null

# --- Row 40: blank spacer row, keeps the hyperlink-column formatting ---
$ws.Range("F30").Copy() | Out-Null
$ws.Range("F40").PasteSpecial(-4122) | Out-Null

# --- Row 41: new component (16-channel analog multiplexer) ---
$ws.Range("D30").Copy() | Out-Null
$ws.Range("D41").PasteSpecial(-4122) | Out-Null
$ws.Range("D41").Value = 1
$ws.Range("E41").Value = 27
$url41 = "https://dratek.cz/arduino/1223-analogovy-multiplexer-16-kanalu-cd74hc4067.html "
$ws.Range("F41").Value = $url41
$ws.Hyperlinks.Add($ws.Range("F41"), $url41.Trim()) | Out-Null
$ws.Range("F30").Copy() | Out-Null
$ws.Range("F41").PasteSpecial(-4122) | Out-Null

# --- Row 42: another BH1750 (reuses the same URL/text as the original component) ---
$ws.Range("D30").Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4122) | Out-Null
$ws.Range("D42").Value = 9
$ws.Range("E42").Value = 67
$url42 = "https://dratek.cz/arduino/902-mereni-intenzity-svetla-bh1750.html"
$ws.Range("F42").Value = $url42
$ws.Hyperlinks.Add($ws.Range("F42"), $url42) | Out-Null
$ws.Range("F30").Copy() | Out-Null
$ws.Range("F42").PasteSpecial(-4122) | Out-Null

# --- Window/selection bookkeeping to mirror the author's last view state ---
$ws.Range("F38").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1

Write-Host "edit complete"
